$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "СОШ №6"
$ws.Range("C2").Value = "13/09/2022"
$ws.Range("D2").Value = "Покупка вещей"
$ws.Range("E2").Value = 1256.832
$ws.Range("F2").Value = 96
$ws.Range("G2").Value = "сраероаеоаеовк"

# Delete row 3 entirely (shifts nothing below it, just removes the row)
$ws.Range("A3:G3").Delete()
